$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 452, shifting existing rows 452:545 down to 453:546.
$ws.Rows(452).Insert()

# Populate the newly inserted row with the new data point.
$ws.Range("A452").Value = 4
$ws.Range("B452").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C452").Value = "Los Lagos"
$ws.Range("D452").Value = 44889
$ws.Range("E452").Value = 10
$ws.Range("F452").Value = 100114001
$ws.Range("G452").Value = "Papa"
$ws.Range("H452").Value = "Rodeo"
$ws.Range("I452").Value = "1a (guarda)"
$ws.Range("J452").Value = 300
$ws.Range("K452").Value = 8000
$ws.Range("L452").Value = 8000
$ws.Range("M452").Value = 8000
$ws.Range("N452").Value = "$/saco 25 kilos"
$ws.Range("O452").Value = "Provincia de Llanquihue"
$ws.Range("P452").Value = 320
$ws.Range("Q452").Value = 25
$ws.Range("R452").Value = "Hortaliza"
